$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell reference -> new text value (Price column D, Volume(1h) column E)
$updates = @{
    "D2" = "30.587.49"
    "E2" = "  +0.38%  "
    "D3" = "2.115.90"
    "E3" = "  +1.10%  "
    "D4" = "1.011"
    "E4" = "  +0.84%  "
    "D5" = "337.46"
    "E5" = "  +2.20%  "
    "D6" = "1.010"
    "E6" = "  +0.84%  "
    "D7" = "0.5250"
    "E7" = "  +0.69%  "
    "D8" = "0.4565"
    "E8" = "  +3.37%  "
    "D9" = "54.63"
    "E9" = "  +1.95%  "
    "D10" = "0.09152"
    "E10" = "  +2.31%  "
    "D11" = "1.175"
    "E11" = "  +1.84%  "
    "D12" = "24.59"
    "E12" = "  +1.23%  "
    "D13" = "2.121.45"
    "E13" = "  +0.81%  "
    "D14" = "6.867"
    "E14" = "  +2.70%  "
    "D15" = "8.110"
    "E15" = "  +5.64%  "
    "E16" = "  +4.85%  "
    "D17" = "97.15"
    "E17" = "  +1.21%  "
    "E18" = "  +0.69%  "
    "D19" = "0.06695"
    "E19" = "  +1.33%  "
    "D20" = "19.47"
    "E20" = "  +1.48%  "
    "D22" = "6.317"
    "E22" = "  +0.92%  "
    "D23" = "30.676.61"
    "E23" = "  +0.51%  "
    "D24" = "12.89"
    "E24" = "  +4.59%  "
    "D25" = "2.360"
    "E25" = "  +1.78%  "
    "D26" = "2.366.61"
    "E26" = "  +0.69%  "
    "D27" = "22.36"
    "E27" = "  +0.45%  "
    "D28" = "164.30"
    "E28" = "  +0.44%  "
    "D29" = "2.545"
    "E29" = "  -0.62%  "
    "D30" = "134.37"
    "E30" = "  +2.11%  "
    "D31" = "1.222"
    "E31" = "  +2.58%  "
    "D32" = "0.1074"
    "E32" = "  +0.37%  "
    "D33" = "1.645"
    "E33" = "  -0.79%  "
    "D34" = "6.379"
    "E34" = "  +3.53%  "
    "D35" = "3.951"
    "E35" = "  +1.40%  "
    "D36" = "10.64"
    "E36" = "  +6.19%  "
    "E37" = "  +7.23%  "
    "D38" = "0.02638"
    "E38" = "  +2.73%  "
    "D39" = "0.06856"
    "E39" = "  +0.22%  "
    "D40" = "0.2332"
    "E40" = "  +3.40%  "
    "D41" = "12.61"
    "E41" = "  +0.28%  "
    "D42" = "0.6878"
    "E42" = "  -0.19%  "
    "E43" = "  +0.46%  "
    "D44" = "14.77"
    "E44" = "  +5.33%  "
    "D45" = "0.6466"
    "E45" = "  +2.01%  "
    "D46" = "2.319"
    "E46" = "  +5.56%  "
    "E47" = "  +21.75%  "
    "E48" = "  +1.78%  "
    "D49" = "1.256"
    "E49" = "  +0.97%  "
    "D50" = "83.61"
    "E50" = "  +2.13%  "
    "D51" = "0.3340"
    "E51" = "  +10.95%  "
}

foreach ($cellRef in $updates.Keys) {
    $range = $ws.Range($cellRef)
    # Force text format so numeric-looking strings (e.g. "1.010") keep their exact
    # textual representation instead of being coerced into a Double.
    $range.NumberFormat = "@"
    $range.Value = $updates[$cellRef]
    $range.NumberFormat = "General"
}

Write-Host "Updated $($updates.Count) cells"
